$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update "Förändrad" (column C) for all existing data rows (2-28) from 45489 to 45490
$ws.Range("C2:C28").Value = 45490

# 2. Add the new data row 29 (A 30076-2024)
$ws.Range("A29").Value = "A 30076-2024"

$ws.Range("B29").Value = 45489
$ws.Range("B29").NumberFormat = "YYYY-MM-DD"

$ws.Range("C29").Value = 45490
$ws.Range("C29").NumberFormat = "YYYY-MM-DD"

$ws.Range("D29").Value = "OKÄNT"
$ws.Range("E29").Value = "OKÄNT"

$ws.Range("G29").Value = 5.6
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = 0
$ws.Range("N29").Value = 0
$ws.Range("O29").Value = 0
$ws.Range("P29").Value = 0
$ws.Range("Q29").Value = 0

$ws.Range("R29").Value = ""
$ws.Range("R29").WrapText = $true

# Row 28 picks up an explicit row height as part of this update
$ws.Rows.Item(28).RowHeight = 15
